$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3
$ws.Range("G3").Value = 1
$ws.Range("G5").Value = 2
$ws.Range("G7").Value = 0
$ws.Range("G8").Value = 1
